# "merubah struktur data barang"
# Insert two new columns (Harga Beli Dus / Harga Beli Pack) between the
# existing "Cash Pack" and "Diskon" columns on the import-barang header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank columns at F:G - this pushes the old F (Diskon) and
# G (Min Qty ...) columns two slots to the right, becoming H and I, and
# keeps their formatting (header style) intact.
$ws.Range("F1:G1").EntireColumn.Insert()

# Fill in the headers for the two newly inserted columns.
$ws.Range("F4").Value2 = "Harga Beli Dus"
$ws.Range("G4").Value2 = "Harga Beli Pack"

# Resize the affected columns (E through H) to fit their new contents,
# matching the width recalculation that happens in Excel after the edit.
$ws.Columns("E:H").AutoFit()

# Reflect the author's final cursor position in the sheet view.
$ws.Range("I7").Select() | Out-Null
